$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 63: A63 and D63 were stored as text ("31"/"11") but should become real numbers ---
$ws.Range("A63").Value = 31
$ws.Range("D63").Value = 11

# --- New data rows 64-70 (new DAP weather readings for 01/19/2025) ---
# Columns: A = Temperature (number), B = Weather Condition (text),
#          C = Date Today (text, formatted like "MM/DD/YYYY"), D = Current Hour (number)

$rows = @(
    @{ Row = 64; A = 31; B = "Mostly Cloudy"; C = "01/19/2025"; D = 12 },
    @{ Row = 65; A = 31; B = "Mostly Cloudy"; C = "01/19/2025"; D = 12 },
    @{ Row = 66; A = 31; B = "Mostly Cloudy"; C = "01/19/2025"; D = 13 },
    @{ Row = 67; A = 32; B = "Mostly Cloudy"; C = "01/19/2025"; D = 13 },
    @{ Row = 68; A = 32; B = "Mostly Cloudy"; C = "01/19/2025"; D = 13 },
    @{ Row = 69; A = 27; B = "Partly Cloudy"; C = "01/19/2025"; D = 19 }
)

foreach ($r in $rows) {
    $ws.Range("A" + $r.Row).Value = $r.A
    $ws.Range("B" + $r.Row).Value = $r.B
    # Keep the date column as literal text, not an auto-converted date serial value.
    $ws.Range("C" + $r.Row).Formula = "'" + $r.C
    $ws.Range("D" + $r.Row).Value = $r.D
}

# Row 70 keeps A and D as plain text (numeric-looking strings), matching the
# trailing raw DAP row that has not been through numeric coercion yet.
$ws.Range("A70").Formula = "'27"
$ws.Range("B70").Value = "Partly Cloudy"
$ws.Range("C70").Formula = "'01/19/2025"
$ws.Range("D70").Formula = "'19"
